$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had two standalone rows (13 and 14) holding just the
# professors' names ("3577649 - Carlos Angelo Nunes" / "1176388 - Luiz Tadeu
# Fernandes Eleno"), sitting between "Docentes responsaveis:" and "Programa
# resumido:". The long paragraph texts for Objetivos / Programa resumido /
# Programa / Bibliografia are dropped, and the short texts that used to
# follow them (dates, names, the "curso sera ministrado..." blurb, the two
# avaliacao paragraphs) slide up to take their place. Net effect: the sheet
# shrinks from 26 to 24 rows.
#
# We first relocate the surviving short text values into their new homes by
# copy/pasting whole cells (so the destination keeps a real "shared string"
# text cell and its original formatting/style -- typing the date-looking
# "01/01/2011" by hand would get auto-converted into a date serial number).
# Moves are ordered bottom of the sheet first so that a later step never
# reads a cell that an earlier step has already overwritten.

# "Para os alunos..." (row 22) -> row 23 (becomes Bibliografia's body later)
$ws.Range("B22").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4163) | Out-Null
$ws.Range("C22").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4163) | Out-Null

# "Serao aplicadas duas avaliacoes..." (row 21) -> row 22
$ws.Range("B21").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4163) | Out-Null
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4163) | Out-Null

# "O curso sera ministrado..." (row 20) -> row 21
$ws.Range("B20").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4163) | Out-Null
$ws.Range("C20").Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4163) | Out-Null

# "1176388 - Luiz Tadeu Fernandes Eleno" (row 14) -> row 20 (under "Metodo:")
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4163) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4163) | Out-Null

# "3577649 - Carlos Angelo Nunes" (row 13) -> row 10 (under "Objetivos:")
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4163) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4163) | Out-Null

# "3577649 - Carlos Angelo Nunes" (row 13) -> row 17 (under "Programa:")
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4163) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4163) | Out-Null

# "01/01/2011" (row 8) -> row 15 (under "Programa resumido:")
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4163) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4163) | Out-Null

# Now drop the two now-redundant standalone rows; everything below shifts up.
$ws.Rows("13:14").Delete()
